$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before the existing "fantasy points" column (G),
# shifting it to column I, so we have room for "height" and "weight".
$ws.Range("G1:H1").EntireColumn.Insert() | Out-Null

# New header labels for the inserted columns.
$ws.Range("G1").Value = "height"
$ws.Range("H1").Value = "weight"

# Copy the style of the header cell (F1) onto the new header cells so they
# match the bold/bordered/centered look of the rest of the header row.
$ws.Range("F1").Copy() | Out-Null
$ws.Range("G1:H1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

$lastRow = 16
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 7).Value = 6.25    # column G: height
    $ws.Cells.Item($r, 8).Value = 210     # column H: weight
}
